$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Find $searchText (literal, not a wildcard pattern) anywhere in the document
# and return a *fresh* Range object covering exactly that match. (Ranges
# returned directly off Document.Content do not support in-place InsertXML
# reliably in this host, so we always re-wrap the found Start/End via
# Document.Range()).
function Find-TextRange($searchText) {
    $scope = $d.Range(0, $d.Content.End)
    $found = $scope.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    return $d.Range($scope.Start, $scope.End)
}

# Replace the exact text $searchText (must be unique in the document) with a
# literal replacement string, keeping it inside whatever single run it was
# already part of (simple 1-for-1 text swap, run count preserved).
function Replace-Text($searchText, $replacement) {
    $scope = $d.Range(0, $d.Content.End)
    $found = $scope.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, $replacement, 2)
    if (-not $found) {
        throw "Text not found for replace: $searchText"
    }
}

# Replace the run(s) spanning $searchText with a brand new sequence of runs
# described by $runsXml (raw OOXML <w:r>...</w:r> fragments concatenated).
# This is used whenever the diff adds/removes run boundaries (splits a
# sentence into extra runs, adds a lastRenderedPageBreak, etc.).
function Replace-WithRuns($searchText, $runsXml) {
    $rng = Find-TextRange $searchText
    $wrapper = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($wrapper)
}

$rPr24 = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr>'
$rPrNoSz = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr>'

# ---------------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------------
Replace-Text "Artificial Intelligence: Navigating the Ethical Conundrum" "The Fascinating Realm of Cells: Microscopic Worlds"

# ---------------------------------------------------------------------------
# Author name: "Kevin Martin" (1 run) -> "Dr" + "." + " Alecia Marshall" (3 runs)
# ---------------------------------------------------------------------------
Replace-WithRuns "Kevin Martin" (
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>Dr</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>.</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve"> Alecia Marshall</w:t></w:r>'
)

# ---------------------------------------------------------------------------
# Email line
# ---------------------------------------------------------------------------
Replace-Text "kevin" "alecia"
Replace-Text "martin@abcxyz" "marshall@educator"
Replace-Text "com" "org"

# ---------------------------------------------------------------------------
# Body paragraph 1
# ---------------------------------------------------------------------------
Replace-Text "As artificial intelligence (AI) strides forward with remarkable progress, its pervasive integration into our lives brings ethical questions to the forefront" "The realm of cells is vast, teeming with countless microscopic worlds that hold the secrets of life"

Replace-Text " The intricate dance between human autonomy and AI-driven decision-making demands careful consideration" " These tiny, intricate structures are the foundation of all living things, from the towering redwood to the minuscule bacterium"

Replace-Text " We must unravel the potential apprehensions and dilemmas while exploring the positive applications of AI that enhance human lives" " Within these minute boundaries, cells perform complex symphonies of biochemical reactions, carrying out functions that sustain life and support growth"

Replace-Text " Striking a balance between progress and responsibility becomes crucial, ensuring AI's impact aligns with ethical standards and societal values" " Exploring the world of cells is a captivating journey into the fundamental mechanisms of biology, revealing the building blocks of life and shedding light on the mysteries of our own existence"

Replace-Text "The transformative power of AI demands a comprehensive examination of its implications" "The diversity of cells is staggering, ranging from simple prokaryotes, like bacteria, to intricate eukaryotes, such as animal and plant cells"

Replace-Text " Our growing dependence on AI-powered systems in domains as diverse as healthcare, finance, and criminal justice mandates ethical scrutiny" " Each cell type is uniquely specialized, adapted to perform specific tasks essential for the survival of the organism"

# This sentence gains two brand-new trailing runs (". " + extra sentence)
Replace-WithRuns " Are AI algorithms biased? Do they perpetuate existing prejudices? What are the consequences of AI-driven decisions gone awry? These are just a few of the ethical landmines we must navigate to ensure AI's ethical integrity" (
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> Specialized cells, such as neurons, facilitate the rapid transmission of information throughout organisms, while muscle cells enable movement and contraction</w:t></w:r>' +
    '<w:r>' + $rPr24 + '<w:t>.</w:t></w:r>' +
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> This exquisite symphony of cells working in concert underscores the intricate complexity of life</w:t></w:r>'
)

Replace-Text "Furthermore, the rise of autonomous AI systems presents unprecedented challenges" "The study of cells has revolutionized our understanding of biology and medicine"

Replace-Text " As these machines become increasingly autonomous, the questions of accountability and liability become tangled" " The development of microscopes has allowed scientists to peer into the inner sanctums of cells, revealing the intricate structures and processes that govern life"

# This sentence also gains two brand-new trailing runs (". " + extra sentence)
Replace-WithRuns " Who bears responsibility when an AI system malfunctions or makes harmful decisions? Legal frameworks and ethical principles must evolve swiftly to address such conundrums" (
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> This knowledge has led to breakthroughs in treating diseases, developing new drugs, and understanding the genetic basis of inheritance</w:t></w:r>' +
    '<w:r>' + $rPr24 + '<w:t>.</w:t></w:r>' +
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> The study of cells continues to unlock mysteries, pushing the boundaries of biological knowledge and offering hope for new treatments and therapies</w:t></w:r>'
)

# ---------------------------------------------------------------------------
# Summary heading - unchanged ("Summary")
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Summary paragraph
# ---------------------------------------------------------------------------
Replace-Text "AI's rapidly expanding role in our lives amplifies the need for ethical considerations" "The microscopic world of cells is a captivating realm of intricate structures and processes that hold the secrets of life"

Replace-Text " We must delve into the complexities of AI-driven decision-making, scrutinize potential biases, and contemplate the consequences of AI's ever-growing autonomy" " From the simplest prokaryotes to the complex eukaryotes, each cell is a finely tuned machine, performing specialized tasks essential for the survival of the organism"

# This sentence is cut short and continues in a brand-new run that also carries <w:lastRenderedPageBreak/>
Replace-WithRuns " This ethical exploration encompasses concerns of privacy, transparency, accountability, and liability" (
    '<w:r>' + $rPrNoSz + '<w:t xml:space="preserve"> The study of cells has revolutionized biology and medicine, leading to groundbreaking </w:t></w:r>' +
    '<w:r>' + $rPrNoSz + '<w:lastRenderedPageBreak/><w:t>discoveries that have improved our understanding of diseases, genetics, and treatments</w:t></w:r>'
)

Replace-Text " Only by grappling with these challenges head-on can we harness AI's potential for progress while safeguarding our values and ensuring its ethical compass remains steadfast" " As we continue to explore the fascinating realm of cells, we unlock the mysteries of life and pave the way for new advancements in healthcare and biological knowledge"

# ---------------------------------------------------------------------------
# New trailing empty paragraph at the end of the document body.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last.Range
$endOfBody = $lastPara.End
$insertPoint = $d.Range($endOfBody, $endOfBody)
$insertPoint.InsertParagraphAfter()
